# Trim the padded whitespace from the header row labels.
# (Market Cap header in H1 already has no surrounding/trailing padding issue
# to fix beyond what's already correct, so it is left as-is.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Stock Code"
$ws.Range("B1").Value = " Company Name"
$ws.Range("C1").Value = " Sector"
$ws.Range("D1").Value = " Open"
$ws.Range("E1").Value = " Close"
$ws.Range("F1").Value = " Volume"
$ws.Range("G1").Value = " Trade Date"

# Move the selection to match where the cursor ended up when the file was saved.
$ws.Range("J9").Select() | Out-Null
